# Adds a new journal entry ("2025. április 9. / Fejlesztői dokumentáció –
# Funkcionális modellek elkészítése / - Kovács Milán") right after the
# existing last entry ("2025. április 6. ...") in the project log.

$d = $word.ActiveDocument

# --- Locate the last (non-empty) paragraph of the main document body. ---
# The Paragraphs collection can contain trailing zero-width "phantom"
# paragraphs (belonging to other stories), so walk back from the end until
# a paragraph with actual content is found.
$n = $d.Paragraphs.Count
$targetIndex = -1
for ($i = $n; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -ne $p.Range.End) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the last content paragraph."
}

$lastPara = $d.Paragraphs.Item($targetIndex)
$lastParaRange = $lastPara.Range

# Sanity-check: this paragraph should be the April 6th entry that currently
# ends the document.
$lastParaText = $lastParaRange.Text
if ($lastParaText -notlike "*Kovács Milán*") {
    throw "Unexpected content in last paragraph: $lastParaText"
}

$rStart = $lastParaRange.Start
$rEnd = $lastParaRange.End

# Range covering the whole paragraph's content, excluding its trailing
# paragraph mark (so the replacement keeps a single, separate paragraph
# mark rather than merging with the next one).
$target = $d.Range($rStart, $rEnd - 1)

# Replace that content with: the original paragraph (re-stated verbatim,
# preserving its rsid attributes) immediately followed by a brand-new
# paragraph for the "2025. április 9." entry. Re-stating the original
# paragraph lets us also relocate the trailing _GoBack bookmark (which
# cannot be reached through the Bookmarks collection because it is
# hidden) from the end of the old last paragraph to the end of the new
# one, exactly like Word would after the user typed the new entry there.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="005D14C8" w:rsidRDefault="005D14C8" w:rsidP="001E2E6C">
<w:pPr><w:keepLines/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr>
<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>2025. április 6.</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:br/><w:t>Fejlesztői dokumentáció &#8211; Követelmény-specifikáció elkészítése</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:br/><w:t>- Kovács Milán</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:keepLines/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr>
<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>2025. április 9.</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:br/><w:t>Fejlesztői dokumentáció &#8211; Funkcionális modellek elkészítése</w:t></w:r>
<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:br/><w:t>- Kovács Milán</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$null = $target.InsertXML($xml)
